$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Initial_Weights (C3) and Final_Weights (E3) change
$ws.Range("C3").Value = "[0.09276893926046914, 0.08168661668268835, 0.04803145896248584]"
$ws.Range("E3").Value = "[5.008104595025083, 5.734199448081903, 10.599637559381387]"

# Row 4: Final_Weights (E4), MSE_Final (G4), RMSE_Final (H4), MAE_Final (I4) change
$ws.Range("E4").Value = "[5.008104595001939, 5.734199448075995, 10.599637559493772]"
$ws.Range("G4").Value = 372.092
$ws.Range("H4").Value = 19.2897
$ws.Range("I4").Value = 15.2952

# Row 5: Initial_Weights (C5), Final_Weights (E5), MSE_Final (G5), RMSE_Final (H5), MAE_Final (I5) change
$ws.Range("C5").Value = "[0.09276893926046914, 0.08168661668268835, 0.04803145896248584]"
$ws.Range("E5").Value = "[5.008104595001945, 5.734199448076011, 10.599637559493772]"
$ws.Range("G5").Value = 372.092
$ws.Range("H5").Value = 19.2897
$ws.Range("I5").Value = 15.2952

# Row 7: Initial_Weights (C7) and Final_Weights (E7) change
$ws.Range("C7").Value = "[0.09276893926046914, 0.08168661668268835, 0.04803145896248584]"
$ws.Range("E7").Value = "[5.008104595020176, 5.734199448068993, 10.599637559494376]"
